# Generated by Katalon AI
# Duplicate the existing data row (row 2) into a new row 3 with the same
# shipping/payment test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "address"
$ws.Range("B3").Value = "Ho Chi Minh"
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "first"
$ws.Range("E3").Value = "last"
$ws.Range("F3").Value = " "

# G3/H3 ("333" / "1111") look numeric, so force them to be entered as text
# (matching the string type used by the rest of the sheet) and then strip
# the number-format override so no stray style is left on the cells.
$ws.Range("G3:H3").NumberFormat = "@"
$ws.Range("G3").Value = "333"
$ws.Range("H3").Value = "1111"
$ws.Range("G3:H3").ClearFormats()
